# Fix Training Data Issue (#48)
#
# The "Date" column (BF) was populated with a string built from the
# source file name ("5-25-2011-12": day-month label + season folder
# "2011-12"), instead of the real ISO game date. NBA.com displays the
# season label right next to the box-score date, and that leaked into
# the scrape. Every data row (BF2:BF31) actually corresponds to the
# same game date, 2012-05-25, so replace the bogus label with the
# correct date, keeping it as plain text (not converting it into a
# date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")

# Writing the literal string via .Formula (wrapped in quotes so it is
# a text-constant, not a date literal) then collapsing the formula
# back down to its value with Copy/PasteSpecial(xlPasteValues) stores
# "2012-05-25" as real text, matching the other BF cells' plain
# (un-formatted) style instead of Excel re-interpreting the assignment
# as a date and reformatting/renumbering the cell.
$dateRange.Formula = "=""2012-05-25"""
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)  # xlPasteValues
